$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.251.51"
$ws.Range("E2").Value = "  +2.79%  "
$ws.Range("D3").Value = "1.893.25"
$ws.Range("E3").Value = "  +0.81%  "
$ws.Range("E4").Value = "  -0.88%  "
$ws.Range("D5").Value = "314.83"
$ws.Range("E5").Value = "  -0.47%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("D8").Value = "0.3906"
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("D9").Value = "0.08421"
$ws.Range("E9").Value = "  +0.14%  "
$ws.Range("E10").Value = "  +1.64%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "6.232"
$ws.Range("E12").Value = "  -0.44%  "
$ws.Range("D13").Value = "1.888.10"
$ws.Range("E13").Value = "  +0.59%  "
$ws.Range("D14").Value = "20.63"
$ws.Range("E14").Value = "  +0.75%  "
$ws.Range("D15").Value = "7.310"
$ws.Range("E15").Value = "  +0.52%  "
$ws.Range("E16").Value = "  -0.90%  "
$ws.Range("D17").Value = "92.89"
$ws.Range("E17").Value = "  +1.89%  "
$ws.Range("D18").Value = "0.00001104"
$ws.Range("E18").Value = "  -0.20%  "
$ws.Range("D19").Value = "0.06738"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "17.81"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("D22").Value = "6.003"
$ws.Range("E22").Value = "  +0.85%  "
$ws.Range("D23").Value = "29.264.22"
$ws.Range("E23").Value = "  +2.56%  "
$ws.Range("E24").Value = "  -0.10%  "
$ws.Range("E25").Value = "  -2.49%  "
$ws.Range("D26").Value = "2.107.10"
$ws.Range("E26").Value = "  +0.67%  "
$ws.Range("D27").Value = "159.30"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "20.85"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").Value = "2.426"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").Value = "127.44"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "1.057"
$ws.Range("E31").Value = "  +0.72%  "
$ws.Range("D32").Value = "0.1046"
$ws.Range("E32").Value = "  -0.66%  "
$ws.Range("D33").Value = "6.175"
$ws.Range("E33").Value = "  +6.67%  "
$ws.Range("E34").Value = "  +0.89%  "
$ws.Range("D35").Value = "0.02475"
$ws.Range("E35").Value = "  +1.70%  "
$ws.Range("D36").Value = "0.06544"
$ws.Range("E36").Value = "  +0.68%  "
$ws.Range("D37").Value = "9.026"
$ws.Range("E37").Value = "  +1.15%  "
$ws.Range("E38").Value = "  -0.10%  "
$ws.Range("D39").Value = "1.223"
$ws.Range("E39").Value = "  +2.78%  "
$ws.Range("D40").Value = "5.134"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("D41").Value = "0.6493"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("D42").Value = "1.228"
$ws.Range("E42").Value = "  -3.20%  "
$ws.Range("D43").Value = "11.24"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("D44").Value = "0.6042"
$ws.Range("E44").Value = "  -0.31%  "
$ws.Range("D45").Value = "13.17"
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "3.675"
$ws.Range("E46").Value = "  -0.98%  "
$ws.Range("D47").Value = "2.043"
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "123.09"
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("D51").Value = "77.30"
$ws.Range("E51").Value = "  +0.38%  "
